$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1029.8572
$ws.Range("I32").Value = 650
$ws.Range("J32").Value = 1181.8
$ws.Range("K32").Value = 650
$ws.Range("L32").Value = 1181.8
$ws.Range("M32").Value = -324
$ws.Range("N32").Value = -1833.8
$ws.Range("H34").Value = 2833.3333
$ws.Range("I34").Value = 2833.3333
$ws.Range("K34").Value = 2833.3333
$ws.Range("M34").Value = -2630.3333
$ws.Range("H36").Value = 2833.3333
$ws.Range("I36").Value = 2833.3333
$ws.Range("K36").Value = 2833.3333
$ws.Range("M36").Value = -2118.3333
$ws.Range("H76").Value = 3300
$ws.Range("I76").Value = 3306.6667
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3306.6667
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2991.6667
$ws.Range("N76").Value = -3830
$ws.Range("H79").Value = 3300
$ws.Range("I79").Value = 3306.6667
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3306.6667
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2214.6667
$ws.Range("N79").Value = -5384
$ws.Range("H98").Value = 967.3333
$ws.Range("I98").Value = 201
$ws.Range("K98").Value = 201
$ws.Range("M98").Value = 1297
$ws.Range("H100").Value = 2914.125
$ws.Range("I100").Value = 2865.2666
$ws.Range("J100").Value = 2995.5557
$ws.Range("K100").Value = 2865.2666
$ws.Range("L100").Value = 2995.5557
$ws.Range("M100").Value = -2324.2666
$ws.Range("N100").Value = -4077.5557
$ws.Range("H112").Value = 4578.041
$ws.Range("J112").Value = 4706.8936
$ws.Range("L112").Value = 14120.6808
$ws.Range("N112").Value = -16336.6808
$ws.Range("H122").Value = 967.3333
$ws.Range("I122").Value = 201
$ws.Range("K122").Value = 603
$ws.Range("M122").Value = 1847
$ws.Range("H138").Value = 2161.0989
$ws.Range("I138").Value = 3497
$ws.Range("J138").Value = 1951.1714
$ws.Range("K138").Value = 10491
$ws.Range("L138").Value = 5853.5142
$ws.Range("M138").Value = -5351
$ws.Range("N138").Value = -16133.5142

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2355.8823
$ws.Range("I102").Value = 2173.077
$ws.Range("K102").Value = 2173.077
$ws.Range("M102").Value = -551.0770000000002
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524
$ws.Range("H113").Value = 40000
$ws.Range("J113").Value = 40000
$ws.Range("L113").Value = 40000
$ws.Range("N113").Value = -48678

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -87
$ws.Range("H8").Value = 6003.6665
$ws.Range("I8").Value = 6003.6665
$ws.Range("K8").Value = 6003.6665
$ws.Range("M8").Value = -5863.6665
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""
$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 4000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 4000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -3761
$ws.Range("N49").Value = ""
$ws.Range("H86").Value = 52633748
$ws.Range("I86").Value = 58825680
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 58825680
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -58824557
$ws.Range("N86").Value = -4546
$ws.Range("H89").Value = 52633748
$ws.Range("I89").Value = 58825680
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 294128400
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -294122784
$ws.Range("N89").Value = -22732
$ws.Range("H99").Value = 1541.4615
$ws.Range("I99").Value = 1566.5555
$ws.Range("K99").Value = 1566.5555
$ws.Range("M99").Value = -68.55549999999994

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5235.125
$ws.Range("I31").Value = 1235.0286
$ws.Range("J31").Value = 11901.952
$ws.Range("K31").Value = 1235.0286
$ws.Range("L31").Value = 11901.952
$ws.Range("M31").Value = -940.0286000000001
$ws.Range("N31").Value = -12491.952
$ws.Range("H34").Value = 5235.125
$ws.Range("I34").Value = 1235.0286
$ws.Range("J34").Value = 11901.952
$ws.Range("K34").Value = 1235.0286
$ws.Range("L34").Value = 11901.952
$ws.Range("M34").Value = -1033.0286
$ws.Range("N34").Value = -12305.952
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 839
$ws.Range("I5").Value = 806.8
$ws.Range("K5").Value = 2420.4
$ws.Range("M5").Value = -2308.4
$ws.Range("H110").Value = 12465.637
$ws.Range("J110").Value = 14012.056
$ws.Range("L110").Value = 42036.16800000001
$ws.Range("N110").Value = -50216.16800000001
$ws.Range("H113").Value = 758.5349
$ws.Range("I113").Value = 469.68182
$ws.Range("J113").Value = 1061.1428
$ws.Range("K113").Value = 1409.04546
$ws.Range("L113").Value = 3183.4284
$ws.Range("M113").Value = 760.95454
$ws.Range("N113").Value = -7523.428400000001
$ws.Range("H122").Value = 8936
$ws.Range("I122").Value = 391.75
$ws.Range("K122").Value = 3525.75
$ws.Range("M122").Value = -1075.75
$ws.Range("H131").Value = 856.13336
$ws.Range("I131").Value = 232.72728
$ws.Range("J131").Value = 1057.8235
$ws.Range("K131").Value = 698.18184
$ws.Range("L131").Value = 3173.4705
$ws.Range("M131").Value = 4341.81816
$ws.Range("N131").Value = -13253.4705
$ws.Range("H134").Value = 3550.8235
$ws.Range("I134").Value = 1901.2174
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5703.6522
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -633.6522000000004
$ws.Range("N134").Value = -31140
$ws.Range("H135").Value = 839
$ws.Range("I135").Value = 806.8
$ws.Range("K135").Value = 7261.2
$ws.Range("M135").Value = -4726.2
$ws.Range("H137").Value = 7758.905
$ws.Range("I137").Value = 11424.363
$ws.Range("J137").Value = 3726.9
$ws.Range("K137").Value = 34273.089
$ws.Range("L137").Value = 11180.7
$ws.Range("M137").Value = -29173.089
$ws.Range("N137").Value = -21380.7
$ws.Range("H140").Value = 1226.7667
$ws.Range("I140").Value = 893.9048
$ws.Range("J140").Value = 2003.4445
$ws.Range("K140").Value = 2681.7144
$ws.Range("L140").Value = 6010.333500000001
$ws.Range("M140").Value = 2498.2856
$ws.Range("N140").Value = -16370.3335

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2540
$ws.Range("I126").Value = 1566.6666
$ws.Range("K126").Value = 4699.9998
$ws.Range("M126").Value = -2229.9998
